$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value reads as a plain number (e.g. "1.00",
# "67.00") are entered with a leading apostrophe so Excel keeps them as
# literal text (matching the source data, which keeps trailing zeros /
# uses "." as a thousands separator). ClearFormats() immediately after
# drops the apostrophe/quote-prefix marker again so no stray style index
# is left on cells that originally carried no explicit style.

$ws.Range("D2").Value = '41.054.45'
$ws.Range("E2").Value = '  -3.98%  '
$ws.Range("D3").Value = '2.454.45'
$ws.Range("E3").Value = '  -3.43%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = "'311.19"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.42%  '
$ws.Range("D6").Value = "'94.03"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -7.30%  '
$ws.Range("D7").Value = "'0.549"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -3.95%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").Value = "'0.501"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -5.37%  '
$ws.Range("D10").Value = "'33.41"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -8.03%  '
$ws.Range("D11").Value = "'0.0780"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.15%  '
$ws.Range("E12").Value = '  -0.62%  '
$ws.Range("D13").Value = "'6.93"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -5.81%  '
$ws.Range("D14").Value = '2.826.26'
$ws.Range("E14").Value = '  -3.60%  '
$ws.Range("D15").Value = '2.444.67'
$ws.Range("E15").Value = '  -5.86%  '
$ws.Range("D16").Value = "'14.46"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -8.21%  '
$ws.Range("D17").Value = "'0.786"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.33%  '
$ws.Range("D18").Value = '41.040.75'
$ws.Range("E18").Value = '  -3.98%  '
$ws.Range("D19").Value = "'6.34"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -6.22%  '
$ws.Range("D20").Value = '0.0₃0914'
$ws.Range("E20").Value = '  -4.09%  '
$ws.Range("D21").Value = "'11.53"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -6.05%  '
$ws.Range("D22").Value = "'67.00"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -3.70%  '
$ws.Range("D23").Value = "'236.20"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -3.55%  '
$ws.Range("D24").Value = "'2.77"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -4.77%  '
$ws.Range("D25").Value = "'1.92"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -6.67%  '
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("D27").Value = "'24.52"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -6.32%  '
$ws.Range("E28").Value = '  -5.19%  '
$ws.Range("D29").Value = "'9.67"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -5.25%  '
$ws.Range("D30").Value = "'36.24"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -7.97%  '
$ws.Range("D31").Value = "'152.75"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.47%  '
$ws.Range("D32").Value = "'5.57"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.18%  '
$ws.Range("E33").Value = '  -1.21%  '
$ws.Range("D34").Value = "'2.53"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -8.64%  '
$ws.Range("D35").Value = "'0.0749"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.82%  '
$ws.Range("D36").Value = "'3.01"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -5.33%  '
$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").Value = "'1.90"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -7.41%  '
$ws.Range("B38").Value = 'Celestia'
$ws.Range("C38").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D38").Value = "'17.09"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -6.53%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").Value = "'0.103"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -8.19%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").Value = "'0.114"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -4.55%  '
$ws.Range("D41").Value = "'4.17"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.55%  '
$ws.Range("D42").Value = "'21.14"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.83%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").Value = '1.964.04'
$ws.Range("E44").Value = '  -1.16%  '
$ws.Range("D45").Value = "'0.0284"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -5.23%  '
$ws.Range("D46").Value = "'3.05"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -8.39%  '
$ws.Range("E47").Value = '  -2.60%  '
$ws.Range("D48").Value = "'76.18"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -5.80%  '
$ws.Range("D49").Value = "'69.48"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -4.55%  '
$ws.Range("D50").Value = "'97.32"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.60%  '
$ws.Range("D51").Value = "'0.179"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -6.82%  '
